$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 740
$ws.Range("C6").Value = 26
$ws.Range("D6").Value = 8.0896686159844062
$ws.Range("A7").Value = 900
$ws.Range("B7").Value = 820
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 10.777087865417259
$ws.Range("A8").Value = 900
$ws.Range("B8").Value = 870
$ws.Range("C8").Value = 32
$ws.Range("D8").Value = 12.29223938056877
$ws.Range("B9").Value = 920
$ws.Range("C9").Value = 27
$ws.Range("D9").Value = 14.02072057161717
$ws.Range("A10").Value = 860
$ws.Range("B10").Value = 980
$ws.Range("C10").Value = 25
$ws.Range("D10").Value = 16.600799007193849
$ws.Range("A11").Value = 850
$ws.Range("B11").Value = 990
$ws.Range("C11").Value = 19
$ws.Range("D11").Value = 17.243623353727081
$ws.Range("A12").Value = 840
$ws.Range("B12").Value = 1000
$ws.Range("C12").Value = 21
$ws.Range("D12").Value = 17.95073013491362
$ws.Range("A13").Value = 820
$ws.Range("B13").Value = 1000
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 19.76891195309544
$ws.Range("A14").Value = 800
$ws.Range("B14").Value = 990
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 23.20901653386435
$ws.Range("A15").Value = 790
$ws.Range("B15").Value = 980
$ws.Range("C15").Value = 15
$ws.Range("D15").Value = 24.2565821356222
$ws.Range("A16").Value = 760
$ws.Range("B16").Value = 930
$ws.Range("C16").Value = 23
$ws.Range("D16").Value = 27.32550418554078
$ws.Range("A17").Value = 730
$ws.Range("B17").Value = 870
$ws.Range("C17").Value = 26
$ws.Range("D17").Value = 30.063546606969091
$ws.Range("A18").Value = 710
$ws.Range("B18").Value = 850
$ws.Range("C18").Value = 23
$ws.Range("D18").Value = 31.218006657885901
$ws.Range("A19").Value = 640
$ws.Range("B19").Value = 830
$ws.Range("C19").Value = 15.17368524526475
$ws.Range("D19").Value = 35.032209998488753
$ws.Range("A20").Value = 620
$ws.Range("C20").Value = 17
$ws.Range("D20").Value = 36.275462044223453
$ws.Range("A21").Value = 570
$ws.Range("B21").Value = 850
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 39.109759311136351
$ws.Range("A22").Value = 520
$ws.Range("B22").Value = 880
$ws.Range("C22").Value = 17
$ws.Range("D22").Value = 42.178681361054927
$ws.Range("A23").Value = 490
$ws.Range("B23").Value = 890
$ws.Range("C23").Value = 17
$ws.Range("D23").Value = 44.038844690565739
$ws.Range("A24").Value = 460
$ws.Range("B24").Value = 890
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 45.974328561533483
$ws.Range("A25").Value = 430
$ws.Range("B25").Value = 870
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 49.109590540197821
$ws.Range("A26").Value = 420
$ws.Range("B26").Value = 850
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 50.972980521447653
$ws.Range("A27").Value = 420
$ws.Range("B27").Value = 840
$ws.Range("C27").Value = 17
$ws.Range("D27").Value = 51.597980521447653
$ws.Range("A28").Value = 420
$ws.Range("B28").Value = 830
$ws.Range("C28").Value = 19
$ws.Range("D28").Value = 52.153536077003203
$ws.Range("B29").Value = 820
$ws.Range("C29").Value = 21
$ws.Range("D29").Value = 52.653536077003203
$ws.Range("B30").Value = 720
$ws.Range("C30").Value = 32
$ws.Range("D30").Value = 56.427120982663581
$ws.Range("A31").Value = 450
$ws.Range("B31").Value = 640
$ws.Range("C31").Value = 22.727959804814699
$ws.Range("D31").Value = 59.549475178138593
$ws.Range("A32").Value = 500
$ws.Range("B32").Value = 580
$ws.Range("C32").Value = 23.889310335394061
$ws.Range("D32").Value = 62.900271847448522
$ws.Range("A33").Value = 540
$ws.Range("B33").Value = 500
$ws.Range("C33").Value = 23
$ws.Range("D33").Value = 66.71532984084692
$ws.Range("A34").Value = 550
$ws.Range("B34").Value = 460
$ws.Range("C34").Value = 20
$ws.Range("D34").Value = 68.63305338764583
$ws.Range("A35").Value = 550
$ws.Range("B35").Value = 450
$ws.Range("C35").Value = 21
$ws.Range("D35").Value = 69.120858265694608
$ws.Range("A36").Value = 540
$ws.Range("B36").Value = 440
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 70.406506958761057
$ws.Range("A37").Value = 520
$ws.Range("B37").Value = 430
$ws.Range("C37").Value = 12
$ws.Range("D37").Value = 73.84661153952996
$ws.Range("A38").Value = 500
$ws.Range("B38").Value = 430
$ws.Range("C38").Value = 16
$ws.Range("D38").Value = 75.27518296810139
$ws.Range("A39").Value = 490
$ws.Range("B39").Value = 430
$ws.Range("C39").Value = 18
$ws.Range("D39").Value = 75.863418262219042
$ws.Range("A40").Value = 480
$ws.Range("B40").Value = 430
$ws.Range("C40").Value = 20
$ws.Range("D40").Value = 76.389734051692727
$ws.Range("A41").Value = 460
$ws.Range("B41").Value = 430
$ws.Range("C41").Value = 23
$ws.Range("D41").Value = 77.319966609832264
$ws.Range("A42").Value = 430
$ws.Range("B42").Value = 430
$ws.Range("C42").Value = 27
$ws.Range("D42").Value = 78.519966609832267
$ws.Range("A43").Value = 390
$ws.Range("B43").Value = 420
$ws.Range("C43").Value = 27
$ws.Range("D43").Value = 80.04704276746844
$ws.Range("A44").Value = 320
$ws.Range("B44").Value = 380
$ws.Range("C44").Value = 25
$ws.Range("D44").Value = 83.147911132198658
$ws.Range("A45").Value = 290
$ws.Range("B45").Value = 340
$ws.Range("C45").Value = 18
$ws.Range("D45").Value = 85.473492527547492
$ws.Range("A46").Value = 280
$ws.Range("B46").Value = 320
$ws.Range("C46").Value = 16
$ws.Range("D46").Value = 86.788826631959139
$ws.Range("A47").Value = 280
$ws.Range("B47").Value = 310
$ws.Range("C47").Value = 11
$ws.Range("D47").Value = 87.529567372699887
$ws.Range("A48").Value = 300
$ws.Range("B48").Value = 270
$ws.Range("C48").Value = 17
$ws.Range("D48").Value = 90.723950197699594
$ws.Range("A49").Value = 310
$ws.Range("B49").Value = 260
$ws.Range("C49").Value = 19
$ws.Range("D49").Value = 91.50962439901798
$ws.Range("A50").Value = 340
$ws.Range("B50").Value = 240
$ws.Range("C50").Value = 24
$ws.Range("D50").Value = 93.186624992257052
$ws.Range("A51").Value = 360
$ws.Range("B51").Value = 230
$ws.Range("C51").Value = 27
$ws.Range("D51").Value = 94.063514395198141
$ws.Range("A52").Value = 400
$ws.Range("B52").Value = 210
$ws.Range("C52").Value = 32
$ws.Range("D52").Value = 95.579492685028512
$ws.Range("A53").Value = 460
$ws.Range("C53").Value = 28
$ws.Range("D53").Value = 97.815560662528299
$ws.Range("A54").Value = 490
$ws.Range("B54").Value = 160
$ws.Range("C54").Value = 24
$ws.Range("D54").Value = 99.202311153091372
$ws.Range("A55").Value = 510
$ws.Range("B55").Value = 150
$ws.Range("C55").Value = 26
$ws.Range("D55").Value = 100.09673834409131
$ws.Range("A56").Value = 550
$ws.Range("B56").Value = 140
$ws.Range("C56").Value = 27
$ws.Range("D56").Value = 101.65262725941869
$ws.Range("A57").Value = 610
$ws.Range("B57").Value = 140
$ws.Range("C57").Value = 28.523339139620909
$ws.Range("D57").Value = 103.8138805604539
$ws.Range("A58").Value = 720
$ws.Range("B58").Value = 150
$ws.Range("C58").Value = 34
$ws.Range("D58").Value = 107.34707669782649
$ws.Range("A59").Value = 780
$ws.Range("B59").Value = 170
$ws.Range("C59").Value = 29
$ws.Range("D59").Value = 109.35487203761591
$ws.Range("A60").Value = 860
$ws.Range("B60").Value = 230
$ws.Range("C60").Value = 24.61337882406934
$ws.Range("D60").Value = 113.0852841919458
$ws.Range("A61").Value = 870
$ws.Range("B61").Value = 240
$ws.Range("C61").Value = 26
$ws.Range("D61").Value = 113.6441141279154
$ws.Range("A65").Value = 940
$ws.Range("B65").Value = 530
$ws.Range("C65").Value = 34
$ws.Range("D65").Value = 124.38804452380521
